# TestPO_09 standard-template migration
# - "Input" sheet: old 17-column layout (거래처명..비고, A:Q) is replaced by the
#   new 16-column standard template layout (발주일자..비고, A:P).
# - "갑지"/"을지" sheets: drop the stray empty I-column cells (I2:I7).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Input" sheet - rebuild with the standard template header + remapped data
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Input")

# Wipe everything first (old header style + the now-dropped columns, e.g. the
# old Q "비고" column) so nothing stale survives the column reshuffle.
$ws.UsedRange.Clear()

$headers = @(
    "발주일자", "납기일자", "거래처명", "거래처 이메일", "납품처명", "납품처 이메일",
    "프로젝트명", "대분류", "중분류", "소분류", "품목명", "규격", "수량", "단가",
    "총금액", "비고"
)
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Row data in new column order:
# 발주일자, 납기일자, 거래처명, 거래처 이메일, 납품처명, 납품처 이메일, 프로젝트명,
# 대분류, 중분류, 소분류, 품목명, 규격, 수량, 단가, 총금액
$rows = @(
    @("2025-09-05", "2025-09-02", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "2. 부자재비", "1) 판넬", "기타", "I-02 120*50*7375*6T 3T", "KS규격-1", 23,  175320, 4435596),
    @("2025-09-12", "2025-09-05", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "2. 부자재비", "1) 판넬", "기타", "3월 절삭",               "KS규격-2", 1,   0,      0),
    @("2025-09-06", "2025-09-16", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "2. 부자재비", "1) 판넬", "기타", "BR-1 B/K 70*74*180*6T", "KS규격-3", 23,  21995,  556473),
    @("2025-09-01", "2025-09-06", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "5. 운반비",   "일반자재", "기타", "4월 운반비",             "KS규격-4", 1,   0,      0),
    @("2025-08-26", "2025-10-01", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "2. 부자재비", "1) 판넬", "기타", "BR-3 B/K 150*90*100*10T","KS규격-5", 46,  3570,   180642),
    @("2025-09-08", "2025-10-12", "센트럴머시너리", "센트럴머시너리@example.com", "힐스테이트 도곡동1차", "delivery@example.com", "힐스테이트 도곡동1차", "5. 운반비",   "일반자재", "기타", "3월 운반비",             "KS규격-6", 1,   0,      0)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $rowData = $rows[$r]

    # Columns A (발주일자) and B (납기일자) hold date-shaped text ("2025-09-05")
    # that must stay literal text, not auto-convert to a date serial - force
    # the cell to Text format before assigning, then drop that formatting
    # again so the cell ends up with no explicit style (matching the
    # template, which has no <c ... s="..."> anywhere). Only touch the two
    # cells we actually (re)formatted so no neighbouring blank cell gets
    # materialized by a wider ClearFormats() call.
    $ws.Cells.Item($rowNum, 1).NumberFormat = "@"
    $ws.Cells.Item($rowNum, 2).NumberFormat = "@"

    for ($col = 1; $col -le $rowData.Length; $col++) {
        $ws.Cells.Item($rowNum, $col).Value = $rowData[$col - 1]
    }

    $ws.Cells.Item($rowNum, 1).ClearFormats()
    $ws.Cells.Item($rowNum, 2).ClearFormats()
}

# The header row inherited no style after UsedRange.Clear(), so it already
# matches the template (no bold/border/centre style on the header cells).

# ---------------------------------------------------------------------------
# 2) "갑지" / "을지" sheets - remove the stray empty inline-string cells in I2:I7
# ---------------------------------------------------------------------------
foreach ($sheetName in @("갑지", "을지")) {
    $sheet = $wb.Worksheets.Item($sheetName)
    $sheet.Range("I2:I7").ClearContents()
}
